$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move "SD card stuff" label up from B22 to B21, and mark the MicroSD
# socket row (row 22) as a "No" (not fitted) item like the rest of the BOM.
$sdCardLabel = $ws.Range("B22").Value()
$ws.Range("B21").Value = $sdCardLabel
$ws.Range("B22").Value = ""
$ws.Range("A22").Value = "N"

# Update the level-shifter part on row 23: it's now the SON/UFDFN-20
# package variant, so update the part code, add the package column, and
# link to its datasheet.
$ws.Range("D23").Value = "TXB0108 (SON)"
$ws.Range("E23").Value = "UFDFN-20"
$ws.Hyperlinks.Add($ws.Range("G23"), "http://www.ti.com/lit/ds/symlink/txb0108.pdf")
$ws.Range("G23").Style = "Hyperlink"

# Reflect the updated scroll position/selection left by the edit.
$ws.Activate()
$ws.Range("D25").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
